# Update handback-status.xlsx datetimes to reflect the regenerated report.
$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for b3fc4203 row (row 3, col G)
$wsOverview.Range("G3").Value = "2016-08-26 08:50:20"

# zh-cn sheet: b3fc4203 row (row 3) - Correspond Handoff Datetime (H) and Correspond Handback DateTime (K)
$wsZhCn.Range("H3").Value = "2016-08-26 08:50:12"
$wsZhCn.Range("K3").Value = "2016-08-26 08:50:42"

# de-de sheet: b3fc4203 row (row 3) - Correspond Handback DateTime (K)
$wsDeDe.Range("K3").Value = "2016-08-26 08:50:49"
